$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.927.09"
$ws.Range("E2").Value = "'  +0.29%  "
$ws.Range("D3").Value = "'1.890.77"
$ws.Range("E3").Value = "'  -0.43%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("D5").Value = "'0.8223"
$ws.Range("E5").Value = "'  +6.12%  "
$ws.Range("D6").Value = "'241.59"
$ws.Range("E6").Value = "'  +0.50%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("D8").Value = "'0.3227"
$ws.Range("E8").Value = "'  +5.29%  "
$ws.Range("D9").Value = "'26.47"
$ws.Range("E9").Value = "'  +3.67%  "
$ws.Range("D10").Value = "'0.07026"
$ws.Range("E10").Value = "'  +2.40%  "
$ws.Range("D11").Value = "'0.08027"
$ws.Range("E11").Value = "'  +0.49%  "
$ws.Range("D12").Value = "'0.7464"
$ws.Range("E12").Value = "'  +1.00%  "
$ws.Range("D13").Value = "'1.891.81"
$ws.Range("E13").Value = "'  -0.70%  "
$ws.Range("D14").Value = "'5.201"
$ws.Range("E14").Value = "'  +0.36%  "
$ws.Range("D15").Value = "'92.04"
$ws.Range("E15").Value = "'  +0.66%  "
$ws.Range("D16").Value = "'29.917.14"
$ws.Range("E16").Value = "'  +0.17%  "
$ws.Range("D17").Value = "'14.03"
$ws.Range("E17").Value = "'  +1.70%  "
$ws.Range("D18").Value = "'5.884"
$ws.Range("E18").Value = "'  -0.33%  "
$ws.Range("D19").Value = "'244.71"
$ws.Range("E19").Value = "'  -0.25%  "
$ws.Range("D20").Value = "'0.000007753"
$ws.Range("E20").Value = "'  +0.38%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "'  +0.03%  "
$ws.Range("D22").Value = "'2.143.04"
$ws.Range("E22").Value = "'  -1.48%  "
$ws.Range("E23").Value = "'  +0.07%  "
$ws.Range("D24").Value = "'6.914"
$ws.Range("E24").Value = "'  -0.40%  "
$ws.Range("D25").Value = "'0.1571"
$ws.Range("E25").Value = "'  +20.72%  "
$ws.Range("D26").Value = "'165.76"
$ws.Range("E26").Value = "'  -0.65%  "
$ws.Range("D27").Value = "'9.184"
$ws.Range("E27").Value = "'  -1.04%  "
$ws.Range("E28").Value = "'  +0.24%  "
$ws.Range("D29").Value = "'2.081"
$ws.Range("E29").Value = "'  +2.50%  "
$ws.Range("D30").Value = "'1.369"
$ws.Range("E30").Value = "'  -1.42%  "
$ws.Range("D31").Value = "'1.513"
$ws.Range("E31").Value = "'  +0.18%  "
$ws.Range("D32").Value = "'4.272"
$ws.Range("D33").Value = "'0.05645"
$ws.Range("E33").Value = "'  +7.34%  "
$ws.Range("D34").Value = "'4.070"
$ws.Range("E34").Value = "'  +0.01%  "
$ws.Range("D35").Value = "'1.271"
$ws.Range("E35").Value = "'  +1.97%  "
$ws.Range("D36").Value = "'0.7301"
$ws.Range("E36").Value = "'  +0.19%  "
$ws.Range("D37").Value = "'2.722"
$ws.Range("E37").Value = "'  -0.23%  "
$ws.Range("D38").Value = "'0.01913"
$ws.Range("E38").Value = "'  +0.12%  "
$ws.Range("D39").Value = "'2.778"
$ws.Range("D40").Value = "'0.4416"
$ws.Range("E40").Value = "'  -0.16%  "
$ws.Range("D41").Value = "'71.83"
$ws.Range("E41").Value = "'  -0.53%  "
$ws.Range("D42").Value = "'5.954"
$ws.Range("E42").Value = "'  -3.89%  "
$ws.Range("D43").Value = "'0.8432"
$ws.Range("E43").Value = "'  +0.59%  "
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "'  +0.00%  "
$ws.Range("D45").Value = "'1.875"
$ws.Range("E45").Value = "'  -0.35%  "
$ws.Range("D48").Value = "'9.694"
$ws.Range("E48").Value = "'  -0.95%  "
$ws.Range("D49").Value = "'989.67"
$ws.Range("E49").Value = "'  +6.28%  "
$ws.Range("D50").Value = "'2.041.72"
$ws.Range("E50").Value = "'  -1.14%  "
$ws.Range("E51").Value = "'  -0.67%  "

# Row 46/47 swap (Quant <-> Aptos)
$ws.Range("B46").Value = "'Aptos"
$ws.Range("C46").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "'7.578"
$ws.Range("E46").Value = "'  -0.32%  "
$ws.Range("B47").Value = "'Quant"
$ws.Range("C47").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "'100.72"
$ws.Range("E47").Value = "'  +0.33%  "
